# Update the company record in row 2 of Hoja1 with the new company's data
# (THE POWER OF LOVE S.A.S -> INVERSIONES IOWA SAS), then leave the
# selection on F2 as it was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A2").Value = "INVERSIONES IOWA SAS"
$ws.Range("B2").Value = 901625453
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "e58a834e-0159-4e8b-97a4-c9fb59bee78a"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 98561334

$ws.Range("F2").Select()
